$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric strings need to be forced to
# Text format first, otherwise Excel auto-converts them to real numbers
# (the source workbook stores every Price/Volume cell as inline text).
$textCells = @("D5", "D6", "D13", "D14", "D20", "D21", "D24", "D26", "D27", "D28", "D33", "D34", "D35", "D38", "D45", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.200.25"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.422.81"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("D5").Value = "563.28"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "144.38"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "2.420.85"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "26.14"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("E15").Value = "  +5.66%  "
$ws.Range("D16").Value = "2.859.12"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "61.908.36"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "2.418.43"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "324.90"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "65.45"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "9.01"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").Value = "587.75"
$ws.Range("E27").Value = "  +14.87%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.527.54"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0944"
$ws.Range("E30").Value = "  +5.50%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +5.58%  "
$ws.Range("D33").Value = "0.150"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  +4.54%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "4.78"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +9.74%  "
$ws.Range("D45").Value = "150.64"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "0.0539"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").Value = "20.38"
$ws.Range("E48").Value = "  +4.71%  "
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "0.0924"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("E51").Value = "  +2.16%  "

# Restore the default cell style on the forced-text cells so they match
# the original (unstyled) formatting instead of keeping an explicit
# "@" text number format.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
